$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Archive"
#
# The localization report had a row ordering / status bug: the handoff
# record for fb45a834-...md and 4879278e-...md were swapped, and the
# 116b6f05-...md / fb45a834-...md rows needed their status flipped from
# "Ready for handoff" to "In Translation" (matching the 116b6f05 handoff
# timestamps). This updates the "Overview" summary sheet plus the two
# per-locale detail sheets (zh-cn, de-de), and refreshes each sheet's
# hyperlinks so the displayed file names follow the corrected rows.
# ---------------------------------------------------------------------------

# ----- Overview sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "fb45a834-f4ec-4e37-b65f-d96cd1b33991.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "2016-03-21 08:34:12"

$ws.Range("A5").Value = "4879278e-d7b4-4e56-be18-3dea002fb98d.md"
$ws.Range("D5").Value = "2016-03-21 08:32:47"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/59e8dcea849139c5d8061b958b60bf6c11ad2f11/e2e/13e5574b-535b-4aa0-8ecf-386518b03605.md", "", "", "13e5574b-535b-4aa0-8ecf-386518b03605.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/33c2b7eea740c7233bb83b079c559b2b10e6b1cc/e2e/116b6f05-847e-41a4-81f3-3850a106a1ca.md", "", "", "116b6f05-847e-41a4-81f3-3850a106a1ca.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d55f623e0d3c3b8c899bf089fc2a786bc437bc7c/e2e/4879278e-d7b4-4e56-be18-3dea002fb98d.md", "", "", "fb45a834-f4ec-4e37-b65f-d96cd1b33991.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/33c2b7eea740c7233bb83b079c559b2b10e6b1cc/e2e/fb45a834-f4ec-4e37-b65f-d96cd1b33991.md", "", "", "4879278e-d7b4-4e56-be18-3dea002fb98d.md") | Out-Null

# ----- zh-cn detail sheet ----------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "fb45a834-f4ec-4e37-b65f-d96cd1b33991.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "fb45a834-f4ec-4e37-b65f-d96cd1b33991.26cdcc324324485fcd3a9a44baa4dfb46ea668e9.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-21 08:34:09"

$ws.Range("A5").Value = "4879278e-d7b4-4e56-be18-3dea002fb98d.md"
$ws.Range("D5").Value = "4879278e-d7b4-4e56-be18-3dea002fb98d.6a44e098b94160a3a7a1c44b0abced56577fcb41.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-21 08:32:44"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/59e8dcea849139c5d8061b958b60bf6c11ad2f11/e2e/13e5574b-535b-4aa0-8ecf-386518b03605.md", "", "", "13e5574b-535b-4aa0-8ecf-386518b03605.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/12710dcc171cc4fe2e068f9a94f5f3c7e0b0bcd3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/13e5574b-535b-4aa0-8ecf-386518b03605.26febb906e2a6ad80bbeb6b7c559113f72c7b118.zh-cn.xlf", "", "", "13e5574b-535b-4aa0-8ecf-386518b03605.26febb906e2a6ad80bbeb6b7c559113f72c7b118.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/cb692da569bd0444e24d6b78bf7117750c26b03f/e2e/13e5574b-535b-4aa0-8ecf-386518b03605.md", "", "", "13e5574b-535b-4aa0-8ecf-386518b03605.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/27a24952e0f568421a45930e0cfacf4d7250e5ef/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/13e5574b-535b-4aa0-8ecf-386518b03605.26febb906e2a6ad80bbeb6b7c559113f72c7b118.zh-cn.xlf", "", "", "13e5574b-535b-4aa0-8ecf-386518b03605.26febb906e2a6ad80bbeb6b7c559113f72c7b118.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/33c2b7eea740c7233bb83b079c559b2b10e6b1cc/e2e/116b6f05-847e-41a4-81f3-3850a106a1ca.md", "", "", "116b6f05-847e-41a4-81f3-3850a106a1ca.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e450584cd4bd1406ab6d473475dda8a46e87c84/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/116b6f05-847e-41a4-81f3-3850a106a1ca.4db4ac9d09437400f767adb0f330a965b176f9b8.zh-cn.xlf", "", "", "116b6f05-847e-41a4-81f3-3850a106a1ca.4db4ac9d09437400f767adb0f330a965b176f9b8.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d55f623e0d3c3b8c899bf089fc2a786bc437bc7c/e2e/4879278e-d7b4-4e56-be18-3dea002fb98d.md", "", "", "fb45a834-f4ec-4e37-b65f-d96cd1b33991.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0220ffb02f777817f7edebd32269bfdf043598b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4879278e-d7b4-4e56-be18-3dea002fb98d.6a44e098b94160a3a7a1c44b0abced56577fcb41.zh-cn.xlf", "", "", "fb45a834-f4ec-4e37-b65f-d96cd1b33991.26cdcc324324485fcd3a9a44baa4dfb46ea668e9.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/33c2b7eea740c7233bb83b079c559b2b10e6b1cc/e2e/fb45a834-f4ec-4e37-b65f-d96cd1b33991.md", "", "", "4879278e-d7b4-4e56-be18-3dea002fb98d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e450584cd4bd1406ab6d473475dda8a46e87c84/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fb45a834-f4ec-4e37-b65f-d96cd1b33991.26cdcc324324485fcd3a9a44baa4dfb46ea668e9.zh-cn.xlf", "", "", "4879278e-d7b4-4e56-be18-3dea002fb98d.6a44e098b94160a3a7a1c44b0abced56577fcb41.zh-cn.xlf") | Out-Null

# ----- de-de detail sheet ----------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "fb45a834-f4ec-4e37-b65f-d96cd1b33991.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "fb45a834-f4ec-4e37-b65f-d96cd1b33991.26cdcc324324485fcd3a9a44baa4dfb46ea668e9.de-de.xlf"
$ws.Range("E4").Value = "2016-03-21 08:34:12"

$ws.Range("A5").Value = "4879278e-d7b4-4e56-be18-3dea002fb98d.md"
$ws.Range("D5").Value = "4879278e-d7b4-4e56-be18-3dea002fb98d.6a44e098b94160a3a7a1c44b0abced56577fcb41.de-de.xlf"
$ws.Range("E5").Value = "2016-03-21 08:32:47"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/59e8dcea849139c5d8061b958b60bf6c11ad2f11/e2e/13e5574b-535b-4aa0-8ecf-386518b03605.md", "", "", "13e5574b-535b-4aa0-8ecf-386518b03605.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a478644f3d3421a9370b0787b3dc3e29f2eae04/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/13e5574b-535b-4aa0-8ecf-386518b03605.26febb906e2a6ad80bbeb6b7c559113f72c7b118.de-de.xlf", "", "", "13e5574b-535b-4aa0-8ecf-386518b03605.26febb906e2a6ad80bbeb6b7c559113f72c7b118.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/20f7a683bb276ad842f0d1a2df3c44916b67f87b/e2e/13e5574b-535b-4aa0-8ecf-386518b03605.md", "", "", "13e5574b-535b-4aa0-8ecf-386518b03605.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ba35d333c20bfff6f4366029d244eea961f4e2df/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/13e5574b-535b-4aa0-8ecf-386518b03605.26febb906e2a6ad80bbeb6b7c559113f72c7b118.de-de.xlf", "", "", "13e5574b-535b-4aa0-8ecf-386518b03605.26febb906e2a6ad80bbeb6b7c559113f72c7b118.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/33c2b7eea740c7233bb83b079c559b2b10e6b1cc/e2e/116b6f05-847e-41a4-81f3-3850a106a1ca.md", "", "", "116b6f05-847e-41a4-81f3-3850a106a1ca.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5045cc51ebcef72093e524fbb9f0880c97ba2c74/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/116b6f05-847e-41a4-81f3-3850a106a1ca.4db4ac9d09437400f767adb0f330a965b176f9b8.de-de.xlf", "", "", "116b6f05-847e-41a4-81f3-3850a106a1ca.4db4ac9d09437400f767adb0f330a965b176f9b8.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d55f623e0d3c3b8c899bf089fc2a786bc437bc7c/e2e/4879278e-d7b4-4e56-be18-3dea002fb98d.md", "", "", "fb45a834-f4ec-4e37-b65f-d96cd1b33991.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d9130f0bddb0729934f2b6fbe438a56589591a7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4879278e-d7b4-4e56-be18-3dea002fb98d.6a44e098b94160a3a7a1c44b0abced56577fcb41.de-de.xlf", "", "", "fb45a834-f4ec-4e37-b65f-d96cd1b33991.26cdcc324324485fcd3a9a44baa4dfb46ea668e9.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/33c2b7eea740c7233bb83b079c559b2b10e6b1cc/e2e/fb45a834-f4ec-4e37-b65f-d96cd1b33991.md", "", "", "4879278e-d7b4-4e56-be18-3dea002fb98d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5045cc51ebcef72093e524fbb9f0880c97ba2c74/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fb45a834-f4ec-4e37-b65f-d96cd1b33991.26cdcc324324485fcd3a9a44baa4dfb46ea668e9.de-de.xlf", "", "", "4879278e-d7b4-4e56-be18-3dea002fb98d.6a44e098b94160a3a7a1c44b0abced56577fcb41.de-de.xlf") | Out-Null

$wb.Worksheets.Item("Overview").Select()
